$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at the top (pushes everything else down by one row),
# so a custom, non-header "Non-heading" row can sit above the real header
# row -- lets downstream logic accept a header-row argument that isn't 0.
$ws.Rows.Item(1).Insert()

$ws.Range("A1").Value = "Non-heading"
$ws.Range("D1").Value = "Ignore"

[void]$ws.Range("A2").Select()
